$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark numeric-looking cells as Text format first so Excel keeps the literal
# string value (matching the inlineStr cells in the source) instead of silently
# converting "243.23", "15", etc. into numeric values.
# NOTE: multi-area ("D2,D3,...") ranges only apply property setters to the
# first area in this runtime, so each cell is formatted individually instead.
$textFormatCells = @(
    "D2",
    "D3",
    "D4",
    "D5",
    "D6",
    "D7",
    "D8",
    "D9",
    "D10",
    "D11",
    "D12",
    "D13",
    "D14",
    "D15",
    "D16",
    "D17",
    "D18",
    "D19",
    "D20",
    "D21",
    "D22",
    "D23",
    "D24",
    "D25",
    "D26",
    "D40",
    "D41",
    "D43",
    "D44",
    "D45",
    "D47",
    "D48"
)
foreach ($cellRef in $textFormatCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}
$ws.Range("G2:G51").NumberFormat = "@"

# Column G (Hora) advances from 14 to 15 for every data row
$ws.Range("G2:G51").Value = "15"

# Per-row coin data updates (price refresh + coin list reshuffle)
$ws.Range("D2").Value = "243.23"
$ws.Range("D3").Value = "23.03"
$ws.Range("D4").Value = "5.429"
$ws.Range("D5").Value = "0.05895"
$ws.Range("D6").Value = "3.450"
$ws.Range("D7").Value = "6.575"
$ws.Range("D8").Value = "0.8112"
$ws.Range("D9").Value = "1.002"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "0.1418"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "0.07440"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "0.03257"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "0.03061"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "0.09338"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Value = "3.854"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "0.001583"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "0.04670"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "0.0005955"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("D19").Value = "0.005860"
$ws.Range("D20").Value = "0.001257"
$ws.Range("D21").Value = "0.004904"
$ws.Range("D22").Value = "0.00006808"
$ws.Range("D23").Value = "3.597"
$ws.Range("D24").Value = "2.134"
$ws.Range("D25").Value = "0.3229"
$ws.Range("D26").Value = "0.1305"
$ws.Range("D40").Value = "0.03943"
$ws.Range("D41").Value = "0.006182"
$ws.Range("D43").Value = "0.002543"
$ws.Range("D44").Value = "0.009148"
$ws.Range("E44").Value = "43LocalTradersLCTBestin24h"
$ws.Range("D45").Value = "0.00005200"
$ws.Range("D47").Value = "0.7506"
$ws.Range("D48").Value = "0.002285"
